$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Locate the first empty row right after the existing tracked-price table.
$usedRows = $ws.UsedRange.Rows.Count
$newRow = $usedRows + 1

$rng = $ws.Range("A$newRow`:D$newRow")

# The sheet stores every column (Date, Price, Discount, Incredible) as text
# (shared strings), e.g. A2="2025-10-16", not a real date. Pre-format the
# new cells as Text so entries like "2026-02-07" aren't auto-converted into
# a date serial number, then write the values and drop back to the sheet's
# normal (unstyled) look, matching the rest of the table.
$rng.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2026-02-07"
$ws.Cells.Item($newRow, 2).Value = "82999000"
$ws.Cells.Item($newRow, 3).Value = "7"
$ws.Cells.Item($newRow, 4).Value = "1"

$rng.Style = "Normal"
